$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-37: price/volume updates only (coin name/link unchanged).
# A leading apostrophe forces the Price column to stay text, matching the
# source data which stores prices as literal strings (e.g. "1.000", "25.958.62").
$ws.Range("D2").Value = "'25.958.62"
$ws.Range("E2").Value = "  -2.48%  "
$ws.Range("D3").Value = "'1.837.07"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Value = "'278.37"
$ws.Range("E5").Value = "  -4.41%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "'0.5091"
$ws.Range("E7").Value = "  -2.98%  "
$ws.Range("D8").Value = "'0.3489"
$ws.Range("E8").Value = "  -5.63%  "
$ws.Range("D9").Value = "'44.79"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").Value = "'0.06819"
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("D11").Value = "'19.91"
$ws.Range("E11").Value = "  -5.67%  "
$ws.Range("D12").Value = "'0.8031"
$ws.Range("E12").Value = "  -8.48%  "
$ws.Range("D13").Value = "'0.07787"
$ws.Range("E13").Value = "  -3.62%  "
$ws.Range("D14").Value = "'1.839.78"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "'5.067"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("D16").Value = "'88.17"
$ws.Range("E16").Value = "  -3.39%  "
$ws.Range("D17").Value = "'0.9994"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "'14.15"
$ws.Range("E18").Value = "  -3.09%  "
$ws.Range("D19").Value = "'0.000008056"
$ws.Range("E19").Value = "  -4.28%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "'26.009.71"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("D22").Value = "'4.775"
$ws.Range("E22").Value = "  -2.98%  "
$ws.Range("D23").Value = "'10.04"
$ws.Range("E23").Value = "  -4.69%  "
$ws.Range("D24").Value = "'6.195"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("D25").Value = "'2.371"
$ws.Range("E25").Value = "  +6.51%  "
$ws.Range("D26").Value = "'143.38"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("D27").Value = "'1.665"
$ws.Range("E27").Value = "  -4.08%  "
$ws.Range("D28").Value = "'17.16"
$ws.Range("E28").Value = "  -3.66%  "
$ws.Range("D29").Value = "'109.59"
$ws.Range("E29").Value = "  -2.92%  "
$ws.Range("D30").Value = "'4.365"
$ws.Range("E30").Value = "  -5.97%  "
$ws.Range("D31").Value = "'4.278"
$ws.Range("E31").Value = "  -6.04%  "
$ws.Range("D32").Value = "'0.08789"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("D33").Value = "'0.04853"
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("D34").Value = "'1.159"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("D35").Value = "'0.7286"
$ws.Range("E35").Value = "  -7.63%  "
$ws.Range("D36").Value = "'2.867"
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("D37").Value = "'3.193"
$ws.Range("E37").Value = "  +0.66%  "

# Rows 38-51: a new "Frax" coin is inserted at row 38, shifting the remaining
# coins down by one row; the former last row (NEARProtocol) drops off the list.
$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").Value = "'0.9989"
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'2.354"
$ws.Range("E39").Value = "  -9.22%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01846"
$ws.Range("E40").Value = "  -3.80%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5131"
$ws.Range("E41").Value = "  -12.93%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.9467"
$ws.Range("E42").Value = "  -9.64%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'117.05"
$ws.Range("E43").Value = "  +2.44%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'6.242"
$ws.Range("E44").Value = "  -2.69%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'7.967"
$ws.Range("E45").Value = "  -6.65%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'0.9994"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4539"
$ws.Range("E47").Value = "  -11.75%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1360"
$ws.Range("E48").Value = "  -7.37%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.245"
$ws.Range("E49").Value = "  -6.42%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'36.12"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05911"
$ws.Range("E51").Value = "  -1.76%  "
